$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the D1 header label from "l27" to "l27.6"
$ws.Range("D1").Value = "l27.6"

# Move the active selection to D1
$ws.Range("D1").Select()
